$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new value in A3 (Bug Number = 1) matching style of the rest of column A
$ws.Range("A3").Value = 1

# Update the active selection on the sheet
$ws.Range("D19").Select()
